$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.536.13"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.902.79"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'238.95"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4909"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'0.2933"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "'0.06700"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "1.885.96"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "'17.06"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("D12").Value = "'0.07337"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "'5.184"
$ws.Range("E13").Value = "  +3.90%  "
$ws.Range("D14").Value = "'88.13"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "'0.6695"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "30.504.70"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'13.51"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "'0.000007891"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'0.9991"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'5.457"
$ws.Range("E20").Value = "  +15.89%  "
$ws.Range("D21").Value = "2.145.92"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'196.45"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'6.134"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "'9.526"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "'162.72"
$ws.Range("E26").Value = "  +3.83%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'1.948"
$ws.Range("E28").Value = "  +6.68%  "
$ws.Range("D29").Value = "'1.473"
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("D31").Value = "'0.09183"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "'4.092"
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").Value = "'0.05177"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").Value = "'0.7451"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").Value = "'2.719"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").Value = "'0.01818"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "'2.070"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "'0.4405"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.939"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'107.03"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").Value = "'69.63"
$ws.Range("E44").Value = "  +22.32%  "
$ws.Range("D45").Value = "'0.9961"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'0.1373"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("D47").Value = "'7.613"
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("D48").Value = "'9.009"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("D49").Value = "'35.14"
$ws.Range("E49").Value = "  +6.25%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'0.3936"
$ws.Range("E51").Value = "  -2.08%  "
